# Update the quarterly unemployment-rate table (Sheet1) from the 2024-Q1
# release (01/04/2024) to the 2024-Q2 release (01/07/2024):
#   - column C (Trimestre) changes from 01/04/2024 to 01/07/2024 for every
#     data row (2-9)
#   - column D (Valor) is refreshed with the new figures
#   - three states swap ranking position, so column A is reordered for
#     rows 4-6 (Distrito Federal / Rio de Janeiro / Rio Grande do Norte)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: re-rank three states (rows 4, 5, 6) ---------------------
$ws.Range("A4").Value = "Rio Grande do Norte"
$ws.Range("A5").Value = "Distrito Federal"
$ws.Range("A6").Value = "Rio de Janeiro"

# --- Column D: new "Valor" figures for the 2024-Q2 release --------------
$ws.Range("D2").Value = 10.5
$ws.Range("D3").Value = 9.699999999999999
$ws.Range("D4").Value = 8.800000000000001
$ws.Range("D5").Value = 8.800000000000001
$ws.Range("D6").Value = 8.5
$ws.Range("D7").Value = 8.4
$ws.Range("D8").Value = 8.699999999999999
$ws.Range("D9").Value = 6.4

# --- Column C: new "Trimestre" date text ---------------------------------
# Force a Text number format first so Excel stores the literal string
# "01/07/2024" instead of auto-converting it to a date serial, then clear
# the formatting override back off so no style change is left behind.
foreach ($r in 2..9) {
    $cell = $ws.Range("C$r")
    $cell.NumberFormat = "@"
    $cell.Value = "01/07/2024"
    $cell.ClearFormats()
}
